$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-09-02T15:43:08-05:00"

# --- "Include from RxNorm" sheet: drop the DB00098 / Thymoglobulin row ---
$rx = $wb.Worksheets.Item("Include from RxNorm")
$rx.Rows.Item(5).Delete()
